# Simulated "time sleeps and screenshots added" edit to the ICGMS test data
# workbook. This mirrors what an interactive automation run through Excel
# would leave behind: different active sheet/cell selections after the run
# pauses (sleeps) at various points, plus one changed input value.

$wb = $excel.ActiveWorkbook

# --- SuperAdmin sheet: selection moves from A2 to A3, tab no longer active ---
$wsSuperAdmin = $wb.Worksheets.Item("SuperAdmin")
$wsSuperAdmin.Range("A3").Select()

# --- Customer sheet: A2 value changed, selection moves to B2 ---
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("A2").Value = "fatbeg106@cream.pink"
$wsCustomer.Range("B4").Value = "0Fb@wAc2"
$wsCustomer.Range("B2").Select()

# --- Repairer sheet: selection moves from A2 to A3 ---
$wsRepairer = $wb.Worksheets.Item("Repairer")
$wsRepairer.Range("A3").Select()

# --- Surveyor sheet: becomes the active sheet/tab ---
$wsSurveyor = $wb.Worksheets.Item("Surveyor")
$wsSurveyor.Activate()
$wsSurveyor.Range("A2").Select()
